$d = $word.ActiveDocument

# --- Simple text replacements (Find & Replace) ---
# 1. Update the activation date
$null = $d.Content.Find.Execute("Ativação: 01/01/2018", $true, $true, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)

# 2. Expand the Objetivos (PT) paragraph with new content about hydrolysates
$null = $d.Content.Find.Execute("Demonstrar as principais etapas no desenvolvimento dos processos bioquímicos industriais abordando aspectos bioquímicos importantes na produção de alimentos, e importantes metabólitos. Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, incluindo as biorefinarias de lignocelulósicos.Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos na resolução de problemas industriais envolvendo processos bioquímicos.", $true, $true, $false, $false, $false, $true, 1, $false, "Demonstrar as principais etapas no desenvolvimento dos processos bioquímicos industriais abordando aspectos bioquímicos importantes na produção de alimentos, e importantes metabólitos. Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, incluindo aspectos bioquímicos de bioprocessos envolvendo a utilização de hidrolisados lignocelulósicos e suas aplicações tecnológicas. Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos na resolução de problemas industriais envolvendo processos bioquímicos.", 2)

# 4. Update Programa resumido (PT)
$null = $d.Content.Find.Execute("Introdução aos processos bioquímicos industriais que incluem o processamento de alimentos, e importantes metabólitos, a manufatura de soros e vacinas, e os conceitos modernos de bioenergia e biorrefinarias.", $true, $true, $false, $false, $false, $true, 1, $false, "Introdução aos processos bioquímicos industriais que incluem o processamento de alimentos, e importantes metabólitos, a manufatura de bioprodutos, e os aspectos bioquímicos de bioprocessos envolvendo bioenergia e biorrefinarias.", 2)

# 5. Update Programa resumido (EN)
$null = $d.Content.Find.Execute("Introduction to industrial biochemical processes which include the food processing and important metabolites, the manufacture of serum and vaccines, the modern concepts of bioenergy and biorefineries.", $true, $true, $false, $false, $false, $true, 1, $false, "Introduction to industrial biochemical processes that include food processing and important metabolites, the manufacture of bioproducts, and the biochemical aspects of bioprocesses involving bioenergy and biorefineries.", 2)

# 6. Rewrite Programa (PT)
$null = $d.Content.Find.Execute("1. Introdução ao Processamento de alimentos: tipos de indústria de alimentos, matériasprimas,fases doprocessamento de produtos alimentícios, conservação/alterações de alimentos, microbiologia dealimentos, alterações bioquímicas em alimentos (oxidação de lipídeos, antioxidantes, escurecimentoenzimático e não enzimático), aflatoxinas, conservantes químicos, toxicantes naturais.2. Discussão e apresentação sobre aspectos bioquímicos importantes na produção de metabólitos por microrganismos e  estudo de casos .3. Manufatura de soros e vacinas Métodosindustriais para a produção de soros e vacinas 4.Biotecnologia de lignocelulósicos: Separação e fermentação das frações e principais processosbioquímicos envolvendo materiais lignocelulósicos.5. Bioenergia, biocombustíveis e biorrefinarias.", $true, $true, $false, $false, $false, $true, 1, $false, "1. Introdução: abordagem geral dos princípios bioquímicos aplicados em diferentes processos2. Processos bioquímicos aplicados à indústria de alimentos: tipos de indústria de alimentos, matéria primas, fases do processamento de produtos alimentícios, conservação/alterações de alimentos 3. Processos bioquímicos nas indústrias de processamento de produtos lácteos, frutas e hortaliças, cacau, produtos gordurosos, produtos desidratados 4. Discussão e apresentação sobre aspectos bioquímicos na produção de bioprodutos: solventes, ácidos orgânicos, fármacos, soros e vacinas, bioinsumos agrícolas e outros de importância industrial5. Bioenergia e biorrefinarias: aspectos bioquímicos de bioprocessos envolvendo a utilização de hidrolisados lignocelulósicos e suas aplicações tecnológicas", 2)

# 7. Rewrite Programa (EN)
$null = $d.Content.Find.Execute("1. Food processing: kinds of food industries, rawmaterials,processing phases of food products,conservation/changes of foods, microbiology of foods,biochemical changes in foods (oxidation oflipids, antioxidants, enzymatic and nonenzymaticdarkening), aflatoxins, chemical preservatives,natural toxicants.2. Discussion and presentation on important biochemical aspects in the production of metabolites by microorganisms and case studies.3. Production of serum and vaccines: industrial methods for production of serum and vaccines –4. Biotechnology of lignocellulose: separation and fermentation of fraction and main biochemicalprocesses involving lignocellulosic materials.5. Bioenergy, biofuels and biorefineries", $true, $true, $false, $false, $false, $true, 1, $false, "1. Introduction: general approach to biochemical principles applied in different processes2. Biochemical processes applied to the food industry: types of food industry, raw materials, stages of food processing, food conservation/modifications3. Biochemical processes in the processing industries of dairy products, fruits and vegetables, cocoa, fatty products, dehydrated products4. Discussion and presentation on biochemical aspects in the production of bioproducts: solvents, organic acids, pharmaceuticals, serums and vaccines, agricultural bioinputs and others of industrial importance5. Bioenergy and biorefineries: biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications", 2)

# 8. Fix missing space before 'sera' in Norma de recuperacao
$null = $d.Content.Find.Execute("A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR)será calculada como MR=(NF+PR)/2", $true, $true, $false, $false, $false, $true, 1, $false, "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF+PR)/2", 2)

# 3. Insert the English translation into the previously-empty italic paragraph
#    that follows the Objetivos (PT) paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim().Length -eq 0 -and $p.Range.Font.Italic) {
        $p.Range.Text = "Demonstrate the main steps in the development of industrial biochemical processes, addressing important biochemical aspects in food production, and important metabolites. Present students with a vision of the potential and strategic applications of modern biotechnology, including biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications. Improve reasoning and awaken students' critical spirit and creativity when solving industrial problems involving biochemical processes"
        break
    }
}

# 9. Replace the Bibliografia paragraph's content: the seven numbered
#    references separated by manual line breaks are replaced with a
#    single run of plain text (three references run together).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("1. BORZANI")) {
        $p.Range.Text = "GAVA, A. J.; SILVA, C. A. B.; FRIAS, J. R. B. Tecnologia de alimentos - princípios e aplicações. São Paulo, Nobel, 2008. ISBN-13: 9788521313823.LIMA, U. A. Biotecnologia Industrial: Processos Fermentativos e Enzimáticos. Volume 3. São Paulo: Editora Edgard Blücher Ltda, 2019. ISBN 9788521214571.Moraes, I. O. Biotecnologia Industrial: Biotecnologia na produção de alimentos. Volume 4. São Paulo: Editora Edgard Blücher Ltda, 2021. ISBN 9786555061529."
        break
    }
}
